$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "06-24-24" (sheet17): only the last saved selection changed.
# ---------------------------------------------------------------------------
$ws24 = $wb.Worksheets.Item("06-24-24")
$ws24.Activate()
$ws24.Range("O20").Select()

# ---------------------------------------------------------------------------
# Sheet "06-25-24" (sheet18): NRFI model retuned - new probabilities/order.
# ---------------------------------------------------------------------------
$ws25 = $wb.Worksheets.Item("06-25-24")
$ws25.Activate()

# Clear out the old Correct/Total (C/D) helper columns before rewriting them,
# since several rows lose values and others gain them after the re-sort.
$ws25.Range("C2:D16").ClearContents()

# New data, sorted by RFPred (column B) descending, after retuning the model.
$ws25.Range("A2").Value = "('CIN', 'PIT')"
$ws25.Range("B2").Value = 0.93200000000000005
$ws25.Range("C2").Value = 0
$ws25.Range("D2").Value = 0

$ws25.Range("A3").Value = "('KC', 'MIA')"
$ws25.Range("B3").Value = 0.92100000000000004
$ws25.Range("C3").Value = 0
$ws25.Range("D3").Value = 0

$ws25.Range("A4").Value = "('AZ', 'MIN')"
$ws25.Range("B4").Value = 0.85399999999999998
$ws25.Range("C4").Value = 0
$ws25.Range("D4").Value = 0

$ws25.Range("A5").Value = "('DET', 'PHI')"
$ws25.Range("B5").Value = 0.83199999999999996
$ws25.Range("C5").Value = 1
$ws25.Range("D5").Value = 1

$ws25.Range("A6").Value = "('CHC', 'SF')"
$ws25.Range("B6").Value = 0.81
$ws25.Range("C6").Value = "N/A"

$ws25.Range("A7").Value = "('BAL', 'CLE')"
$ws25.Range("B7").Value = 0.748
$ws25.Range("C7").Value = 0
$ws25.Range("D7").Value = 0

$ws25.Range("A8").Value = "('LAA', 'OAK')"
$ws25.Range("B8").Value = 0.68

$ws25.Range("A9").Value = "('ATL', 'STL')"
$ws25.Range("B9").Value = 0.63800000000000001

$ws25.Range("A10").Value = "('MIL', 'TEX')"
$ws25.Range("B10").Value = 0.60799999999999998

$ws25.Range("A11").Value = "('SEA', 'TB')"
$ws25.Range("B11").Value = 0.58299999999999996

$ws25.Range("A12").Value = "('CWS', 'LAD')"
$ws25.Range("B12").Value = 0.53700000000000003

$ws25.Range("A13").Value = "('BOS', 'TOR')"
$ws25.Range("B13").Value = 0.47

$ws25.Range("A14").Value = "('SD', 'WSH')"
$ws25.Range("B14").Value = 0.41199999999999998

$ws25.Range("A15").Value = "('COL', 'HOU')"
$ws25.Range("B15").Value = 0.222
$ws25.Range("C15").Value = 0
$ws25.Range("D15").Value = 1

$ws25.Range("A16").Value = "('NYM', 'NYY')"
$ws25.Range("B16").Value = 0.099000000000000005
$ws25.Range("C16").Value = 0
$ws25.Range("D16").Value = 1

# Recompute the summary formulas that depend on column D.
$ws25.Range("E2").Formula = "=COUNTIF(D:D, 1)"
$ws25.Range("F2").Formula = "=COUNTA(D2:D20)"
$ws25.Range("H2").Formula = "=E2/F2*100"

$excel.CalculateFull()

$ws25.Range("E4").Select()
